$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.672.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.714.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.13%  "
$ws.Range("E4").Value = "  -0.89%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.521"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("E7").Value = "  -1.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.71"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +10.65%  "
$ws.Range("E9").Value = "  +3.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0626"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0894"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.945.83"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.707.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.24"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.569"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.642.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "240.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0747"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.992"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.60"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.114"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("E29").Value = "  -0.71%  "
$ws.Range("E30").Value = "  +1.40%  "
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.41"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.546.71"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.28"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.69"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.969"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.615"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.41%  "
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0174"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("E40").Value = "  +1.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "71.18"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.68%  "
$ws.Range("E42").Value = "  +5.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.992"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.850.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.790"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "90.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.77%  "
$ws.Range("E50").Value = "  +2.61%  "
$ws.Range("E51").Value = "  -4.97%  "
